$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header cell A1: "Citation" -> "en"
$ws.Range("A1").Value = "en"

# Clear out the leftover empty formatted rows below the data (rows 21, 25, 26, 30, 31)
# so the sheet's used range shrinks back down to A1:D20.
$ws.Rows("21:31").Delete()

# Scroll the view down and select A16, matching the saved view state.
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Range("A16").Select()
